$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 (pushes existing rows 18-42 down to 19-43),
# mirroring the new weekly price record added ahead of the existing ones.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44897
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = "Espárragos"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 1200
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("N18").Value = "$/kilo"
$ws.Range("O18").Value = "Región de Ñuble"
$ws.Range("P18").Value = 950
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
